$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Edit 1 (slide 10): fix the typo "Unknow " -> "Unknown " in the content
# placeholder's 5th paragraph ("Unknow people service will log the time ...").
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange
$para5 = $tr10.Paragraphs(5, 1)
$typo = $tr10.Characters($para5.Start, 7)
if ($typo.Text -eq "Unknow ") {
    $typo.Text = "Unknown "
}

# ---------------------------------------------------------------------------
# Edit 2 (slide 7): append a hyperlink to the Wikipedia article on software
# design patterns after the "Design patterns in code level" bullet.
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(2)
$tr7 = $sh7.TextFrame.TextRange
$para6 = $tr7.Paragraphs(6, 1)

$origStart = $para6.Start
$origVisibleLen = $para6.Text.TrimEnd([char]13).Length

$urlText = "https://en.wikipedia.org/wiki/Software_design_pattern"
$null = $para6.InsertAfter("    " + $urlText)

$full7 = $sh7.TextFrame.TextRange
$spacesStart = $origStart + $origVisibleLen
$spacesRange = $full7.Characters($spacesStart, 4)

$urlStart = $spacesStart + 4
$urlRange = $full7.Characters($urlStart, $urlText.Length)

$action = $urlRange.ActionSettings(1)
$action.Hyperlink.Address = $urlText
$urlRange.Font.Underline = $true
